$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 0.42801434716386477
$ws.Range("C2").Value = 2.0273389758469595
$ws.Range("D2").Value = 0.38689077983612485
$ws.Range("E2").Value = 1.3939160069157488

$ws.Range("B3").Value = 0.91898135611724596
$ws.Range("C3").Value = 2.6799012164395939
$ws.Range("D3").Value = 1.1116130973260185
$ws.Range("E3").Value = 1.0633454225808727

$ws.Range("B1:E3").Select()
